# Auto-generated edit script updating market price / profit columns (H-N)
# across multiple job sheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1460.579
$ws.Range("I5").Value = 1603.4706
$ws.Range("K5").Value = 1603.4706
$ws.Range("M5").Value = -1488.4706
$ws.Range("H9").Value = 208.33333
$ws.Range("I9").Value = 53.444443
$ws.Range("K9").Value = 53.444443
$ws.Range("M9").Value = 115.555557
$ws.Range("H74").Value = 10967.777
$ws.Range("I74").Value = 11142.235
$ws.Range("K74").Value = 11142.235
$ws.Range("M74").Value = -10206.235
$ws.Range("H77").Value = 10967.777
$ws.Range("I77").Value = 11142.235
$ws.Range("K77").Value = 55711.175
$ws.Range("M77").Value = -51031.175
$ws.Range("H80").Value = 630.9091
$ws.Range("I80").Value = 508.5
$ws.Range("J80").Value = 700.8570999999999
$ws.Range("K80").Value = 1525.5
$ws.Range("L80").Value = 2102.5713
$ws.Range("M80").Value = -527.5
$ws.Range("N80").Value = -4098.5713
$ws.Range("H83").Value = 630.9091
$ws.Range("I83").Value = 508.5
$ws.Range("J83").Value = 700.8570999999999
$ws.Range("K83").Value = 4576.5
$ws.Range("L83").Value = 6307.7139
$ws.Range("M83").Value = 415.5
$ws.Range("N83").Value = -16291.7139
$ws.Range("H98").Value = 1861.0869
$ws.Range("I98").Value = 1147.7894
$ws.Range("J98").Value = 5249.25
$ws.Range("K98").Value = 1147.7894
$ws.Range("L98").Value = 5249.25
$ws.Range("M98").Value = 350.2106000000001
$ws.Range("N98").Value = -8245.25
$ws.Range("H113").Value = 3012.7
$ws.Range("I113").Value = 2564.8333
$ws.Range("J113").Value = 3684.5
$ws.Range("K113").Value = 2564.8333
$ws.Range("L113").Value = 3684.5
$ws.Range("M113").Value = 689.1667000000002
$ws.Range("N113").Value = -10192.5
$ws.Range("H122").Value = 1861.0869
$ws.Range("I122").Value = 1147.7894
$ws.Range("J122").Value = 5249.25
$ws.Range("K122").Value = 3443.3682
$ws.Range("L122").Value = 15747.75
$ws.Range("M122").Value = -993.3681999999999
$ws.Range("N122").Value = -20647.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1366.9048
$ws.Range("I2").Value = 1614.7693
$ws.Range("K2").Value = 1614.7693
$ws.Range("M2").Value = -1501.7693
$ws.Range("H14").Value = 4006
$ws.Range("I14").Value = 4006
$ws.Range("K14").Value = 4006
$ws.Range("M14").Value = -3831
$ws.Range("H116").Value = 1366.9048
$ws.Range("I116").Value = 1614.7693
$ws.Range("K116").Value = 1614.7693
$ws.Range("M116").Value = 679.2307000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1366.9048
$ws.Range("I3").Value = 1614.7693
$ws.Range("K3").Value = 1614.7693
$ws.Range("M3").Value = -1500.7693
$ws.Range("H46").Value = 35000
$ws.Range("J46").Value = 35000
$ws.Range("L46").Value = 35000
$ws.Range("N46").Value = -35596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6332.6665
$ws.Range("J86").Value = 6332.6665
$ws.Range("L86").Value = 6332.6665
$ws.Range("N86").Value = -8578.666499999999
$ws.Range("H89").Value = 6332.6665
$ws.Range("J89").Value = 6332.6665
$ws.Range("L89").Value = 31663.3325
$ws.Range("N89").Value = -42895.3325
$ws.Range("H105").Value = 1562.2
$ws.Range("I105").Value = 1569.1111
$ws.Range("K105").Value = 1569.1111
$ws.Range("M105").Value = 177.8888999999999
$ws.Range("H134").Value = 5419.68
$ws.Range("I134").Value = 4247.316
$ws.Range("K134").Value = 12741.948
$ws.Range("M134").Value = -10206.948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 427
$ws.Range("I5").Value = 385.77777
$ws.Range("K5").Value = 1157.33331
$ws.Range("M5").Value = -1045.33331
$ws.Range("H68").Value = 3037.4736
$ws.Range("I68").Value = 8750
$ws.Range("J68").Value = 2365.4119
$ws.Range("K68").Value = 26250
$ws.Range("L68").Value = 7096.2357
$ws.Range("M68").Value = -25439
$ws.Range("N68").Value = -8718.235700000001
$ws.Range("H71").Value = 3037.4736
$ws.Range("I71").Value = 8750
$ws.Range("J71").Value = 2365.4119
$ws.Range("K71").Value = 78750
$ws.Range("L71").Value = 21288.7071
$ws.Range("M71").Value = -74694
$ws.Range("N71").Value = -29400.7071
$ws.Range("H132").Value = 2921.5625
$ws.Range("I132").Value = 2606.2
$ws.Range("J132").Value = 3447.1667
$ws.Range("K132").Value = 23455.8
$ws.Range("L132").Value = 31024.5003
$ws.Range("M132").Value = -20925.8
$ws.Range("N132").Value = -36084.5003
$ws.Range("H135").Value = 427
$ws.Range("I135").Value = 385.77777
$ws.Range("K135").Value = 3471.99993
$ws.Range("M135").Value = -936.9999299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4168.4546
$ws.Range("I80").Value = 3512.5557
$ws.Range("K80").Value = 3512.5557
$ws.Range("M80").Value = -2514.5557
$ws.Range("H83").Value = 4168.4546
$ws.Range("I83").Value = 3512.5557
$ws.Range("K83").Value = 17562.7785
$ws.Range("M83").Value = -12570.7785

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47386.543
$ws.Range("I7").Value = 60419.445
$ws.Range("J7").Value = 8287.833000000001
$ws.Range("K7").Value = 60419.445
$ws.Range("L7").Value = 8287.833000000001
$ws.Range("M7").Value = -60307.445
$ws.Range("N7").Value = -8511.833000000001
$ws.Range("H82").Value = 2298.9
$ws.Range("I82").Value = 5835.6
$ws.Range("K82").Value = 5835.6
$ws.Range("M82").Value = -5474.6
$ws.Range("H85").Value = 2298.9
$ws.Range("I85").Value = 5835.6
$ws.Range("K85").Value = 5835.6
$ws.Range("M85").Value = -4587.6
$ws.Range("H126").Value = 47386.543
$ws.Range("I126").Value = 60419.445
$ws.Range("J126").Value = 8287.833000000001
$ws.Range("K126").Value = 181258.335
$ws.Range("L126").Value = 24863.499
$ws.Range("M126").Value = -178788.335
$ws.Range("N126").Value = -29803.499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 35000
$ws.Range("I51").Value = 20000
$ws.Range("J51").Value = 50000
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 50000
$ws.Range("M51").Value = -19490
$ws.Range("N51").Value = -51020
$ws.Range("H100").Value = 836
$ws.Range("I100").Value = 642
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1284
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -743
$ws.Range("N100").Value = -5082
